$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.130.45"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "1.608.39"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  -0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.72"
$ws.Range("E5").Value = "  +1.63%  "

$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.482"
$ws.Range("E7").Value = "  +1.01%  "

$ws.Range("E8").Value = "  +2.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0619"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.41"
$ws.Range("E10").Value = "  +3.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "1.831.65"
$ws.Range("E12").Value = "  +0.56%  "

$ws.Range("D13").Value = "1.602.53"
$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").Value = "26.142.12"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.86"
$ws.Range("E17").Value = "  +0.40%  "

$ws.Range("E18").Value = "  +2.53%  "

$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "198.66"
$ws.Range("E20").Value = "  +5.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  +2.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.49"
$ws.Range("E22").Value = "  +1.97%  "

$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("E24").Value = "  +3.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.79"
$ws.Range("E25").Value = "  +0.84%  "

$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("E27").Value = "  -0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.22"
$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.51"
$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("E31").Value = "  +2.28%  "

$ws.Range("E32").Value = "  +2.71%  "

$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("E34").Value = "  +4.29%  "

$ws.Range("E35").Value = "  -2.47%  "

$ws.Range("D36").Value = "1.105.99"
$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("E38").Value = "  +1.62%  "

$ws.Range("E39").Value = "  +2.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("E41").Value = "  -1.19%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").Value = "  +8.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.17"
$ws.Range("E43").Value = "  +2.04%  "

$ws.Range("D44").Value = "1.743.35"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.15"
$ws.Range("E45").Value = "  -2.53%  "

$ws.Range("E46").Value = "  -4.34%  "

$ws.Range("E47").Value = "  +6.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.93"
$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.408"
$ws.Range("E50").Value = "  -0.52%  "

$ws.Range("E51").Value = "  -0.64%  "
